$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18/19 - timeparse function: swap set ordering and scalar value
$ws.Range("E18").Value = "{'num', 'empty'}"
$ws.Range("E19").Value = "num"

# Rows 27/28 - line_splitter function: swap set ordering and scalar value
$ws.Range("E27").Value = "{'any', 'Tuple[None]'}"
$ws.Range("E28").Value = "any"

# Rows 31/32 - json_splitter function: swap set ordering and scalar value
$ws.Range("E31").Value = "{'any', 'Tuple[None]'}"
$ws.Range("E32").Value = "any"

# Rows 39/40 - parse_seconds_float function: swap set ordering and scalar value
$ws.Range("E39").Value = "{'num', 'empty'}"
$ws.Range("E40").Value = "num"

# Row 211 - move "Scalpel Accuracy:" label and its value two columns to the right,
# and fix the accuracy value (was incorrectly 511.76, now correctly 83.65000000000001)
$ws.Range("C211").Value = $null
$ws.Range("D211").Value = $null
$ws.Range("E211").Value = "Scalpel Accuracy:"
$ws.Range("F211").Value = 83.65000000000001

# Row 212 - fix label wording
$ws.Range("E212").Value = "Accuracy vs PyType"
